# Insert a new weekly price record at row 60 for "Macroferia Regional de
# Talca - Berenjena". Existing rows 60:63 shift down to 61:64.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(60).Insert()

$ws.Range("A60").Value = 5
$ws.Range("B60").Value = "Macroferia Regional de Talca"
$ws.Range("C60").Value = "Maule"
$ws.Range("D60").Value = 44461
$ws.Range("E60").Value = 7
$ws.Range("F60").Value = 100112001
$ws.Range("G60").Value = "Berenjena"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 200
$ws.Range("K60").Value = 7000
$ws.Range("L60").Value = 7000
$ws.Range("M60").Value = 7000
$ws.Range("N60").Value = "`$/caja 50 unidades"
$ws.Range("O60").Value = "Región de Arica y Parinacota"
$ws.Range("P60").Value = 140
$ws.Range("Q60").Value = 50
$ws.Range("R60").Value = "Hortaliza"
